$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-02-23 Sunday"; new = "2025-02-24 Monday"},
    @{old = "575×8="; new = "994×4="},
    @{old = "254×3="; new = "900×6="},
    @{old = "353×4="; new = "310×6="},
    @{old = "870×8="; new = "636×2="},
    @{old = "175×6="; new = "802×6="},
    @{old = "593×6="; new = "994×7="},
    @{old = "307×2="; new = "363×4="},
    @{old = "653×7="; new = "912×3="},
    @{old = "862×8="; new = "392×6="},
    @{old = "579×4="; new = "438×4="},
    @{old = "603×3="; new = "243×2="},
    @{old = "506×6="; new = "764×5="},
    @{old = "550×7="; new = "861×9="},
    @{old = "845×2="; new = "215×5="},
    @{old = "456×4="; new = "285×3="},
    @{old = "780×2="; new = "317×5="},
    @{old = "444×9="; new = "835×2="},
    @{old = "333×3="; new = "957×5="},
    @{old = "621×8="; new = "610×3="},
    @{old = "180×4="; new = "594×4="},
    @{old = "759×5="; new = "357×9="},
    @{old = "392×5="; new = "336×3="},
    @{old = "947×8="; new = "603×8="},
    @{old = "487×2="; new = "855×2="},
    @{old = "863×6="; new = "396×6="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}
